# Generate Report for Handoff
# Updates the localization-status workbook after a new handoff pass for
# b.md: the file is now "Ready for handoff" (instead of already handed
# back), a new (b.*) handoff xliff exists for both target languages, and
# an error was detected because the existing handback file is stale.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/b342cccd51e418fd787d5d107c1c5f8de858cc03/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/8e8ac05a49cc49b5406ed6289b59b0839751135d/e2e/b.md."

# --- Overview sheet: b.md row (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-03 02:32:52"

# --- zh-cn sheet: b.md row (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe keeps "False" stored as text (not auto-coerced to a
# Boolean); re-applying the Normal style afterwards drops the resulting
# quote-prefix formatting so the cell style stays the default.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-03 02:32:42"
$wsZhCn.Range("P3").Value = $errorDetail
# 39.17 characters round-trips (via the engine's px-quantised column-width
# conversion) to a raw OOXML column width of exactly 40.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: b.md row (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-03 02:32:52"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
